$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Nota dos alunos"
$ws.Range("B3").Value = "Professor"
$ws.Range("C3").Value = "Descrição"
